# campaign/VolunteerInformation.xlsx
# "removing badly named poll phone files"
#
# - Row 53 (Laraib): the contact number cell gains a second phone number,
#   turning it from a plain number into a text value.
# - A new row 55 is appended for "Faraz", a phone canvasser, with his
#   e-mail (hyperlinked) and phone number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update E53: was a single numeric phone number, becomes text with two numbers ---
$ws.Cells.Item(53, 5).Value = "6477195000 / 9054059413"

# --- Add the new row 55 ---
$ws.Cells.Item(55, 1).Value = "Faraz"
$ws.Cells.Item(55, 2).Value = "faraz_tahir@hotmail.com"
$ws.Cells.Item(55, 4).Value = "Phone Canvasser?"
$ws.Cells.Item(55, 5).Value = 6475206843

# Hyperlink the new e-mail address (mirrors how other e-mails in the sheet are linked)
$ws.Hyperlinks.Add($ws.Cells.Item(55, 2), "mailto:faraz_tahir@hotmail.com") | Out-Null

# Re-apply the cell formatting used by the neighbouring rows so the new
# row matches the look of the rest of the table (the hyperlink creation
# above resets the font on B55, so this must happen afterwards).
$ws.Cells.Item(54, 4).Copy()
$ws.Cells.Item(55, 4).PasteSpecial(-4122)

$ws.Cells.Item(53, 2).Copy()
$ws.Cells.Item(55, 2).PasteSpecial(-4122)

# Match the row height used by the other data rows in this block
$ws.Rows.Item(55).RowHeight = 18.75

# Update the selection/scroll position to reflect where the editor left off
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("E55").Select() | Out-Null
